$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:A19").Select()
$ws.Columns("A:A").Insert()

$ws.Range("A1").Value = "Match ID"
$ws.Range("A1").Font.Bold = $true

$ws.Range("A2").Font.Bold = $true
$ws.Range("A3").Font.Bold = $true

$ws.Range("A4:A19").Value = 13
$ws.Range("A4:A19").Font.Bold = $true

$ws.Range("A20").Value = 13
$ws.Rows(20).AutoFit()

$ws.Range("A1:A19").Select()
